$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; everything currently at row 32 and below
# shifts down by one (matches the diff: dimension grows from T122 to T123).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record. Columns A, B, C,
# E-J hold the same constant values (market/product identifiers) found in
# every other row of this subset sheet.
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44708
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100104
$ws.Range("H32").Value = "Frutos de pepita"
$ws.Range("I32").Value = 100104001
$ws.Range("J32").Value = "Granada"
$ws.Range("K32").Value = "Wonderfull"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 60
$ws.Range("N32").Value = 19000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 19417
$ws.Range("Q32").Value = "$/bandeja 15 kilos granel"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 1294
$ws.Range("T32").Value = 15
